$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 ("Population Median Age" / mage): recategorize from
# "Population Change" (demo) to "Population Stock" (phcs), matching the
# adjacent rows 9 and 11, and flip the "period" flag off.
$ws.Range("K10").Value = 0
$ws.Range("N10").Value = "phcs"
$ws.Range("O10").Value = "Population & Human Capital Stocks"
$ws.Range("P10").Value = "Population Stock"

# --- Row 31 ("Net Migration" / net): now flagged as available by age too.
$ws.Range("E31").Value = 1

# --- Row 26 ("Mean Age at Childbearing" / macb): updated definition text,
# data now available in all scenarios (not just by level of educational
# attainment).
$ws.Range("Q26").Value = "The mean age of mothers at the birth of their children observed in a five-year period. Available in all scenarios and at country level. It is expressed in years."

# --- Update last-selected cell to reflect where the author left off.
$ws.Range("A10").Select()
